$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 19 (shifts RealCash..TicketList down by 2)
$ws.Rows("19:20").Insert()

# Fill the two newly inserted rows: StarCandy / AccStarCandy (DOUBLE), mirroring Gold/AccGold pattern
$ws.Range("A19").Value = "StarCandy"
$ws.Range("B19").Value = "DOUBLE"

$ws.Range("A20").Value = "AccStarCandy"
$ws.Range("B20").Value = "DOUBLE"
$ws.Range("E20").Value = "Model"

# Rename KingdomObj -> KingdomItem (this row was row 23, now shifted to row 25)
$ws.Range("A25").Value = "KingdomItemList"
$ws.Range("B25").Value = "LIST:KingdomItemPacket"

# Match the final cursor/selection position recorded in the saved workbook
[void]$ws.Range("E17").Select()
